# Generate Report for Handoff
# Updates status text from "In Translation" to "Ready for handoff" and
# refreshes the handoff timestamps, then widens the affected "Status"
# columns to fit the new (longer) text.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Status text: "In Translation" -> "Ready for handoff" ---
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$zhcn.Range("C2").Value = "Ready for handoff"
$dede.Range("C2").Value = "Ready for handoff"

# --- Refresh handoff timestamps ---
$overview.Range("G2").Value = "2016-08-22 10:19:38"
$dede.Range("H2").Value = "2016-08-22 10:19:38"
$zhcn.Range("H2").Value = "2016-08-22 10:19:33"

# --- Widen the Status columns so the longer text fits ---
$overview.Columns.Item(5).ColumnWidth = 17.2159881591797
$overview.Columns.Item(6).ColumnWidth = 17.2159881591797
$zhcn.Columns.Item(3).ColumnWidth = 17.2159881591797
$dede.Columns.Item(3).ColumnWidth = 17.2159881591797
